# Auto-generated edit script: updates Leve profit-calculator values
# across all 8 Disciple of the Hand crafting sheets, refreshed by the
# scheduled market-data runner (current Universalis prices changed).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1940693.6
$ws.Range("J17").Value = 1986881.5
$ws.Range("L17").Value = 5960644.5
$ws.Range("N17").Value = -5960980.5

$ws.Range("H33").Value = 462.26666
$ws.Range("I33").Value = 471.31033
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 471.31033
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -242.31033
$ws.Range("N33").Value = -658

$ws.Range("H107").Value = 345.83334
$ws.Range("I107").Value = 226.8125
$ws.Range("J107").Value = 583.875
$ws.Range("K107").Value = 226.8125
$ws.Range("L107").Value = 583.875
$ws.Range("M107").Value = 1693.1875
$ws.Range("N107").Value = -4423.875

$ws.Range("H132").Value = 3863.75
$ws.Range("I132").Value = 3805.3726
$ws.Range("J132").Value = 4092.7693
$ws.Range("K132").Value = 11416.1178
$ws.Range("L132").Value = 12278.3079
$ws.Range("M132").Value = -8886.1178
$ws.Range("N132").Value = -17338.3079

$ws.Range("H135").Value = 1517.3928
$ws.Range("I135").Value = 1054.0476
$ws.Range("J135").Value = 2907.4285
$ws.Range("K135").Value = 9486.428400000001
$ws.Range("L135").Value = 26166.8565
$ws.Range("M135").Value = -6951.428400000001
$ws.Range("N135").Value = -31236.8565

$ws.Range("H137").Value = 1107.5405
$ws.Range("I137").Value = 851.84
$ws.Range("J137").Value = 1640.25
$ws.Range("K137").Value = 2555.52
$ws.Range("L137").Value = 4920.75
$ws.Range("M137").Value = -5.519999999999982
$ws.Range("N137").Value = -10020.75

$ws.Range("H141").Value = 3005.4285
$ws.Range("I141").Value = 1643.2727
$ws.Range("J141").Value = 8000
$ws.Range("K141").Value = 4929.8181
$ws.Range("L141").Value = 24000
$ws.Range("M141").Value = 250.1818999999996
$ws.Range("N141").Value = -34360


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1467.96
$ws.Range("I61").Value = 1454.3043
$ws.Range("J61").Value = 1625
$ws.Range("K61").Value = 1454.3043
$ws.Range("L61").Value = 1625
$ws.Range("M61").Value = -1242.3043
$ws.Range("N61").Value = -2049

$ws.Range("H132").Value = 3031.923
$ws.Range("I132").Value = 2862.2173
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 8586.651899999999
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -6056.651899999999
$ws.Range("N132").Value = -18059

$ws.Range("H136").Value = 1467.96
$ws.Range("I136").Value = 1454.3043
$ws.Range("J136").Value = 1625
$ws.Range("K136").Value = 4362.9129
$ws.Range("L136").Value = 4875
$ws.Range("M136").Value = -1812.9129
$ws.Range("N136").Value = -9975


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3126769.2
$ws.Range("I86").Value = 4083431.5
$ws.Range("J86").Value = 1672.4667
$ws.Range("K86").Value = 4083431.5
$ws.Range("L86").Value = 1672.4667
$ws.Range("M86").Value = -4082308.5
$ws.Range("N86").Value = -3918.4667

$ws.Range("H89").Value = 3126769.2
$ws.Range("I89").Value = 4083431.5
$ws.Range("J89").Value = 1672.4667
$ws.Range("K89").Value = 20417157.5
$ws.Range("L89").Value = 8362.333500000001
$ws.Range("M89").Value = -20411541.5
$ws.Range("N89").Value = -19594.3335

$ws.Range("H94").Value = 5473.864
$ws.Range("I94").Value = 951.05554
$ws.Range("K94").Value = 951.05554
$ws.Range("M94").Value = -500.05554


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2095.9575
$ws.Range("I31").Value = 1218.3125
$ws.Range("J31").Value = 3968.2666
$ws.Range("K31").Value = 1218.3125
$ws.Range("L31").Value = 3968.2666
$ws.Range("M31").Value = -923.3125
$ws.Range("N31").Value = -4558.2666

$ws.Range("H34").Value = 2095.9575
$ws.Range("I34").Value = 1218.3125
$ws.Range("J34").Value = 3968.2666
$ws.Range("K34").Value = 1218.3125
$ws.Range("L34").Value = 3968.2666
$ws.Range("M34").Value = -1016.3125
$ws.Range("N34").Value = -4372.2666

$ws.Range("H132").Value = 1405.8572
$ws.Range("I132").Value = 1226.2
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3678.6
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1148.6
$ws.Range("N132").Value = -20057

$ws.Range("H134").Value = 2106.64
$ws.Range("I134").Value = 2119.0205
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 6357.0615
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -3822.0615
$ws.Range("N134").Value = -9570


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 717.4
$ws.Range("I5").Value = 545.7778
$ws.Range("J5").Value = 1232.2667
$ws.Range("K5").Value = 1637.3334
$ws.Range("L5").Value = 3696.800099999999
$ws.Range("M5").Value = -1525.3334
$ws.Range("N5").Value = -3920.800099999999

$ws.Range("H131").Value = 940.75
$ws.Range("I131").Value = 569.53845
$ws.Range("J131").Value = 1150.5652
$ws.Range("K131").Value = 1708.61535
$ws.Range("L131").Value = 3451.6956
$ws.Range("M131").Value = 3331.38465
$ws.Range("N131").Value = -13531.6956

$ws.Range("H135").Value = 717.4
$ws.Range("I135").Value = 545.7778
$ws.Range("J135").Value = 1232.2667
$ws.Range("K135").Value = 4912.000199999999
$ws.Range("L135").Value = 11090.4003
$ws.Range("M135").Value = -2377.000199999999
$ws.Range("N135").Value = -16160.4003


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 48792344
$ws.Range("I70").Value = 103679320
$ws.Range("J70").Value = 3922.2222
$ws.Range("K70").Value = 103679320
$ws.Range("L70").Value = 3922.2222
$ws.Range("M70").Value = -103679050
$ws.Range("N70").Value = -4462.2222

$ws.Range("H73").Value = 48792344
$ws.Range("I73").Value = 103679320
$ws.Range("J73").Value = 3922.2222
$ws.Range("K73").Value = 103679320
$ws.Range("L73").Value = 3922.2222
$ws.Range("M73").Value = -103678384
$ws.Range("N73").Value = -5794.2222

$ws.Range("H80").Value = 3299.1667
$ws.Range("I80").Value = 3638.8462
$ws.Range("J80").Value = 2416
$ws.Range("K80").Value = 3638.8462
$ws.Range("L80").Value = 2416
$ws.Range("M80").Value = -2640.8462
$ws.Range("N80").Value = -4412

$ws.Range("H83").Value = 3299.1667
$ws.Range("I83").Value = 3638.8462
$ws.Range("J83").Value = 2416
$ws.Range("K83").Value = 18194.231
$ws.Range("L83").Value = 12080
$ws.Range("M83").Value = -13202.231
$ws.Range("N83").Value = -22064

$ws.Range("H126").Value = 2486.111
$ws.Range("I126").Value = 1918.75
$ws.Range("J126").Value = 2940
$ws.Range("K126").Value = 5756.25
$ws.Range("L126").Value = 8820
$ws.Range("M126").Value = -3286.25
$ws.Range("N126").Value = -13760


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13165918
$ws.Range("I132").Value = 15633700
$ws.Range("J132").Value = 4416.3335
$ws.Range("K132").Value = 46901100
$ws.Range("L132").Value = 13249.0005
$ws.Range("M132").Value = -46898570
$ws.Range("N132").Value = -18309.0005

$ws.Range("H136").Value = 8024.6816
$ws.Range("I136").Value = 9930.286
$ws.Range("J136").Value = 4689.875
$ws.Range("K136").Value = 29790.858
$ws.Range("L136").Value = 14069.625
$ws.Range("M136").Value = -27240.858
$ws.Range("N136").Value = -19169.625


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2049.16
$ws.Range("I132").Value = 1040.6945
$ws.Range("J132").Value = 4642.357
$ws.Range("K132").Value = 3122.0835
$ws.Range("L132").Value = 13927.071
$ws.Range("M132").Value = -592.0835000000002
$ws.Range("N132").Value = -18987.071

